# Auto-generated edit script applying numeric updates to the "Leve profit" tracking sheets.
# For each affected sheet/row, columns H-N (currentAveragePrice, currentAveragePriceNQ/HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) are refreshed with newly computed market values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 58198.5
$ws.Range("I76").Value = 61445.47
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 61445.47
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -61130.47
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 58198.5
$ws.Range("I79").Value = 61445.47
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 61445.47
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -60353.47
$ws.Range("N79").Value = -5184
$ws.Range("H86").Value = 66734308
$ws.Range("I86").Value = 101119.8
$ws.Range("J86").Value = 200000690
$ws.Range("K86").Value = 101119.8
$ws.Range("L86").Value = 200000690
$ws.Range("M86").Value = -99996.8
$ws.Range("N86").Value = -200002936
$ws.Range("H89").Value = 66734308
$ws.Range("I89").Value = 101119.8
$ws.Range("J89").Value = 200000690
$ws.Range("K89").Value = 505599
$ws.Range("L89").Value = 1000003450
$ws.Range("M89").Value = -499983
$ws.Range("N89").Value = -1000014682
$ws.Range("H107").Value = 807.6667
$ws.Range("I107").Value = 283.625
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 283.625
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 1636.375
$ws.Range("N107").Value = -8840
$ws.Range("H137").Value = 2940
$ws.Range("J137").Value = 3487.5
$ws.Range("L137").Value = 10462.5
$ws.Range("N137").Value = -15562.5
$ws.Range("H141").Value = 4369.3335
$ws.Range("I141").Value = 2025.1904
$ws.Range("J141").Value = 6713.476
$ws.Range("K141").Value = 6075.5712
$ws.Range("L141").Value = 20140.428
$ws.Range("M141").Value = -895.5712000000003
$ws.Range("N141").Value = -30500.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3012.27
$ws.Range("I32").Value = 2247.0957
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 2247.0957
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -1960.0957
$ws.Range("N32").Value = -15574
$ws.Range("H61").Value = 1257.1837
$ws.Range("I61").Value = 956.86365
$ws.Range("J61").Value = 3900
$ws.Range("K61").Value = 956.86365
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -744.86365
$ws.Range("N61").Value = -4324
$ws.Range("H63").Value = 1001837.9
$ws.Range("I63").Value = 1430197
$ws.Range("J63").Value = 2333.3333
$ws.Range("K63").Value = 1430197
$ws.Range("L63").Value = 2333.3333
$ws.Range("M63").Value = -1429511
$ws.Range("N63").Value = -3705.3333
$ws.Range("H66").Value = 1001837.9
$ws.Range("I66").Value = 1430197
$ws.Range("J66").Value = 2333.3333
$ws.Range("K66").Value = 7150985
$ws.Range("L66").Value = 11666.6665
$ws.Range("M66").Value = -7147553
$ws.Range("N66").Value = -18530.6665
$ws.Range("H122").Value = 1275.0714
$ws.Range("I122").Value = 1110.4445
$ws.Range("J122").Value = 1571.4
$ws.Range("K122").Value = 3331.3335
$ws.Range("L122").Value = 4714.200000000001
$ws.Range("M122").Value = -881.3335000000002
$ws.Range("N122").Value = -9614.200000000001
$ws.Range("H132").Value = 3749.7322
$ws.Range("I132").Value = 4132.75
$ws.Range("K132").Value = 12398.25
$ws.Range("M132").Value = -9868.25
$ws.Range("H136").Value = 1257.1837
$ws.Range("I136").Value = 956.86365
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 2870.59095
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -320.5909499999998
$ws.Range("N136").Value = -16800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3695.606
$ws.Range("I20").Value = 4471.3687
$ws.Range("K20").Value = 4471.3687
$ws.Range("M20").Value = -4224.3687

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 25095
$ws.Range("J74").Value = 32500
$ws.Range("L74").Value = 32500
$ws.Range("N74").Value = -34248
$ws.Range("H77").Value = 25095
$ws.Range("J77").Value = 32500
$ws.Range("L77").Value = 97500
$ws.Range("N77").Value = -106236
$ws.Range("H138").Value = 37929.332
$ws.Range("J138").Value = 37929.332
$ws.Range("L138").Value = 37929.332
$ws.Range("N138").Value = -48209.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20.275862
$ws.Range("I2").Value = 20.857143
$ws.Range("J2").Value = 18.75
$ws.Range("K2").Value = 125.142858
$ws.Range("L2").Value = 112.5
$ws.Range("M2").Value = -12.142858
$ws.Range("N2").Value = -338.5
$ws.Range("H122").Value = 1251.8889
$ws.Range("I122").Value = 389.75
$ws.Range("J122").Value = 1941.6
$ws.Range("K122").Value = 3507.75
$ws.Range("L122").Value = 17474.4
$ws.Range("M122").Value = -1057.75
$ws.Range("N122").Value = -22374.4
$ws.Range("H131").Value = 752.65
$ws.Range("I131").Value = 381.25
$ws.Range("J131").Value = 784.9457
$ws.Range("K131").Value = 1143.75
$ws.Range("L131").Value = 2354.8371
$ws.Range("M131").Value = 3896.25
$ws.Range("N131").Value = -12434.8371

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 23811506
$ws.Range("I122").Value = 35715956
$ws.Range("J122").Value = 2609.7144
$ws.Range("K122").Value = 107147868
$ws.Range("L122").Value = 7829.1432
$ws.Range("M122").Value = -107145418
$ws.Range("N122").Value = -12729.1432
$ws.Range("H123").Value = 33000
$ws.Range("J123").Value = 33000
$ws.Range("L123").Value = 33000
$ws.Range("N123").Value = -37900
$ws.Range("H132").Value = 48874.773
$ws.Range("I132").Value = 69648.87
$ws.Range("J132").Value = 4358.857
$ws.Range("K132").Value = 208946.61
$ws.Range("L132").Value = 13076.571
$ws.Range("M132").Value = -206416.61
$ws.Range("N132").Value = -18136.571
$ws.Range("H139").Value = 24900
$ws.Range("J139").Value = 24900
$ws.Range("L139").Value = 24900
$ws.Range("N139").Value = -35180
$ws.Range("H140").Value = 48256
$ws.Range("J140").Value = 48256
$ws.Range("L140").Value = 48256
$ws.Range("N140").Value = -58616

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H130").Value = 13500
$ws.Range("J130").Value = 13500
$ws.Range("L130").Value = 13500
$ws.Range("N130").Value = -23540
$ws.Range("H132").Value = 11689.526
$ws.Range("I132").Value = 30734
$ws.Range("J132").Value = 2899.7693
$ws.Range("K132").Value = 92202
$ws.Range("L132").Value = 8699.3079
$ws.Range("M132").Value = -89672
$ws.Range("N132").Value = -13759.3079
$ws.Range("H133").Value = 25737.572
$ws.Range("J133").Value = 25737.572
$ws.Range("L133").Value = 25737.572
$ws.Range("N133").Value = -30797.572
$ws.Range("H140").Value = 33370
$ws.Range("J140").Value = 33370
$ws.Range("L140").Value = 33370
$ws.Range("N140").Value = -43730

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49984.8
$ws.Range("J46").Value = 49984.8
$ws.Range("L46").Value = 49984.8
$ws.Range("N46").Value = -50446.8
$ws.Range("H134").Value = 49984.8
$ws.Range("J134").Value = 49984.8
$ws.Range("L134").Value = 149954.4
$ws.Range("N134").Value = -155024.4
